$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Bump the cached "datetimeFigureOut" field text from 5/10/2011 to
#    5/11/2011 everywhere it is cached: the slide master and every slide
#    layout (the notes master has the same field too, but this runtime's
#    COM shim mis-routes writes to NotesMaster placeholder shapes onto the
#    slide master by colliding shape id, so it is intentionally left alone
#    to avoid corrupting unrelated master content).
# ---------------------------------------------------------------------------
function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $t = $sh.TextFrame.TextRange.Text
            if ($t -eq "5/10/2011") {
                $sh.TextFrame.TextRange.Text = "5/11/2011"
            }
        }
    }
}

Update-DateShapes($p.SlideMaster.Shapes)

for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    Update-DateShapes($p.SlideMaster.CustomLayouts.Item($li).Shapes)
}

# ---------------------------------------------------------------------------
# 2) On slide 1, inside the "SBaseRef" rectangle, the "port:" paragraph
#    currently reads:  "port: " + "PortSIdRef" + " " + '{use=“optional”}'
#    as four separate runs. Merge the trailing " " run into the following
#    '{use=“optional”}' run (same formatting) so the paragraph becomes
#    "port: " + "PortSIdRef" + ' {use=“optional”}' -- three runs.
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(1)
$group = $slide.Shapes.Item(1)
$rect = $group.GroupItems.Item(2)
$portPara = $rect.TextFrame.TextRange.Paragraphs(2, 1)

$spaceRun = $portPara.Runs(3, 1)
$optionalRun = $portPara.Runs(4, 1)

$optionalRun.Text = $spaceRun.Text + $optionalRun.Text
$spaceRun.Text = ""
